$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the description text column (F) for rows 7-13: "observation point" -> "observation area"
$ws.Range("F7").Value = "Anomaly or deviation from the 1981 - 2010 average for sea surface temperature taken at the Niño 1+2  observation area (0-10°South)(90°West-80°West) "
$ws.Range("F8").Value = "Sea surface temperature in degrees Celsius taken at the Niño 3 observation area (5°North-5°South)(150°West-90°West)"
$ws.Range("F9").Value = "Anomaly or deviation from the 1981 - 2010 average for sea surface temperature taken at the Niño 3 observation area (5°North-5°South)(150°West-90°West)"
$ws.Range("F10").Value = "Sea surface temperature in degrees Celsius taken at the Niño 3.4 observation area (5°North-5°South)(170-120°West"
$ws.Range("F11").Value = "Anomaly or deviation from the 1981 - 2010 average for sea surface temperature taken at the Niño 3.4 observation area (5°North-5°South)(170-120°West"
$ws.Range("F12").Value = "Sea surface temperature in degrees Celsius taken at the Niño 4 observation area (5°North-5°South)"
$ws.Range("F13").Value = "Anomaly or deviation from the 1981 - 2010 average for sea surface temperature taken at the Niño 4 observation area (5°North-5°South)"

# Update the selected cell in the sheet view from I16 to H12
$ws.Range("H12").Select()
